# Refresh the cryptocurrency price/volume snapshot (values scraped
# on 2024-03-25) and re-sort two coin pairs whose ranking flipped
# between runs (OKB/Bittensor at rows 33-34, THORChain/LidoDAOToken
# at rows 49-50).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several "Price" cells are digit strings Excel would otherwise
# auto-convert to numbers (dropping trailing/leading zeros, e.g.
# "9.48" or "0.0000280"). Write them quote-prefixed to force text,
# then restore the Normal style so no quote-prefix formatting
# lingers on the cell (matches the source file, where these are
# plain unstyled text cells).
function Set-TextCell($ref, $text) {
    $ws.Range($ref).Value = "'" + $text
    $ws.Range($ref).Style = 'Normal'
}

$ws.Range('D2').Value = '67.389.07'
$ws.Range('E2').Value = '  +4.71%  '
$ws.Range('D3').Value = '3.487.89'
$ws.Range('E3').Value = '  +4.65%  '
$ws.Range('E4').Value = '  +0.11%  '
Set-TextCell 'D5' '586.54'
$ws.Range('E5').Value = '  +5.92%  '
Set-TextCell 'D6' '186.37'
$ws.Range('E6').Value = '  +7.26%  '
Set-TextCell 'D7' '0.635'
$ws.Range('E7').Value = '  +1.50%  '
$ws.Range('D8').Value = '3.483.03'
$ws.Range('E8').Value = '  +4.80%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('E10').Value = '  +1.26%  '
Set-TextCell 'D11' '0.654'
$ws.Range('E11').Value = '  +3.13%  '
Set-TextCell 'D12' '56.44'
$ws.Range('E12').Value = '  +5.64%  '
Set-TextCell 'D13' '0.0000280'
$ws.Range('E13').Value = '  +0.52%  '
Set-TextCell 'D14' '9.48'
$ws.Range('E14').Value = '  +4.19%  '
$ws.Range('D15').Value = '4.047.43'
$ws.Range('E15').Value = '  +4.82%  '
Set-TextCell 'D16' '18.97'
$ws.Range('E16').Value = '  +4.39%  '
$ws.Range('D17').Value = '3.488.07'
$ws.Range('E17').Value = '  +4.34%  '
$ws.Range('D18').Value = '67.473.33'
$ws.Range('E18').Value = '  +4.81%  '
Set-TextCell 'D19' '12.22'
$ws.Range('E19').Value = '  +3.77%  '
$ws.Range('E20').Value = '  -0.92%  '
$ws.Range('E21').Value = '  +3.49%  '
Set-TextCell 'D22' '491.33'
$ws.Range('E22').Value = '  +8.43%  '
Set-TextCell 'D23' '5.38'
$ws.Range('E23').Value = '  +7.38%  '
Set-TextCell 'D24' '16.68'
$ws.Range('E24').Value = '  +20.00%  '
Set-TextCell 'D25' '4.43'
$ws.Range('E25').Value = '  +9.02%  '
Set-TextCell 'D26' '90.33'
$ws.Range('E26').Value = '  +2.67%  '
$ws.Range('E27').Value = '  +2.32%  '
Set-TextCell 'D28' '11.02'
$ws.Range('E28').Value = '  +4.13%  '
Set-TextCell 'D29' '9.19'
$ws.Range('E29').Value = '  +6.75%  '
Set-TextCell 'D30' '31.54'
$ws.Range('E30').Value = '  +1.58%  '
$ws.Range('E31').Value = '  +10.43%  '
Set-TextCell 'D32' '11.77'
$ws.Range('E32').Value = '  +2.98%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell 'D33' '597.20'
$ws.Range('E33').Value = '  +4.91%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell 'D34' '64.22'
$ws.Range('E34').Value = '  +3.77%  '
$ws.Range('E35').Value = '  +4.65%  '
$ws.Range('E36').Value = '  +6.17%  '
$ws.Range('E37').Value = '  -0.06%  '
Set-TextCell 'D38' '36.74'
$ws.Range('E38').Value = '  +3.61%  '
Set-TextCell 'D39' '3.56'
$ws.Range('E39').Value = '  +1.41%  '
$ws.Range('E40').Value = '  +5.63%  '
$ws.Range('D41').Value = '0.0₃0768'
$ws.Range('E41').Value = '  +5.12%  '
$ws.Range('D42').Value = '3.262.49'
$ws.Range('E42').Value = '  +6.29%  '
Set-TextCell 'D43' '2.93'
$ws.Range('E43').Value = '  +6.26%  '
Set-TextCell 'D44' '0.0431'
$ws.Range('E44').Value = '  +3.83%  '
$ws.Range('E45').Value = '  +3.43%  '
Set-TextCell 'D46' '2.79'
$ws.Range('E46').Value = '  +22.84%  '
Set-TextCell 'D47' '3.27'
$ws.Range('E47').Value = '  +2.07%  '
$ws.Range('E48').Value = '  +1.34%  '
$ws.Range('B49').Value = 'LidoDAOToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 'D49' '3.27'
$ws.Range('E49').Value = '  +13.14%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextCell 'D50' '8.80'
$ws.Range('E50').Value = '  +7.84%  '
$ws.Range('E51').Value = '  +0.12%  '
